# Apply the "Major update. Option for parallel computing" edit to Params.xlsx
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("30 Minutes")

# --- Row 2: update simulation parameters ---
$ws1.Range("C2").Value = -8
$ws1.Range("D2").Value = 4
$ws1.Range("E2").Value = -18
$ws1.Range("H2").ClearFormats()
$ws1.Range("H2").Value = 0.2
$ws1.Range("I2").Formula = "=PI()/24"
$ws1.Range("J2").Value = 3
$ws1.Range("L2").ClearContents()

# --- Rows 3-9: remove stale duplicated parameter columns (B:K), keep A/L/M/N and the P:AB block ---
$ws1.Range("B3:K9").Clear()
$ws1.Range("L3:M9").ClearContents()

# --- Column H (dt) no longer used below row 2; drop the leftover styled blanks entirely ---
$ws1.Range("H3:H25").Clear()

# --- A couple of stray styled-but-empty K cells get fully removed too ---
$ws1.Range("K13").Clear()
$ws1.Range("K14").Clear()
$ws1.Range("K19").Clear()

# --- Row 10 no longer carries a duplicated parameter block or the old "Changed Thermal Model" note ---
$ws1.Range("P10:Y10").Clear()
$ws1.Range("AB10").Clear()

# --- Re-select "30 Minutes" as the active sheet/cell ---
$ws1.Activate()
$ws1.Range("B2").Select()

# --- Add the new (empty) "Log" sheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$logSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$logSheet.Name = "Log"

# Keep "30 Minutes" as the active/selected sheet after adding "Log"
$ws1.Activate()
